$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 0.0257
$ws.Range("H4").Value = 0.257

$ws.Range("G6").Value = 0.098
$ws.Range("H6").Value = 0.098

$ws.Range("G7").Value = 0.0621
$ws.Range("H7").Value = 0.1863

$ws.Range("G9").Value = 0.094
$ws.Range("H9").Value = 1.13

$ws.Range("G11").Value = 0.1059
$ws.Range("H11").Value = 0.1059

$ws.Range("G24").Value = 0.7432
$ws.Range("H24").Value = 0.7432

$ws.Range("G27").Value = 0.4349
$ws.Range("H27").Value = 0.8698

$ws.Range("G30").Value = 0.7586
$ws.Range("H30").Value = 0.7586

$ws.Range("G32").Value = 0.4762
$ws.Range("H32").Value = 0.9524

$ws.Range("G35").Value = 0.0087
$ws.Range("H35").Value = 0.087

$ws.Range("G36").Value = 0.0689
$ws.Range("H36").Value = 0.0689

$ws.Range("G37").Value = 0.0929
$ws.Range("H37").Value = 0.0929

$ws.Range("G39").Value = 0.01
$ws.Range("H39").Value = 0.1

$ws.Range("G40").Value = 0.0955
$ws.Range("H40").Value = 0.0955

$ws.Range("G41").Value = 0.012
$ws.Range("H41").Value = 0.12

$ws.Range("G42").Value = 0.012
$ws.Range("H42").Value = 0.12

$ws.Range("G43").Value = 0.061
$ws.Range("H43").Value = 0.061

$ws.Range("G47").Value = 0.4461
$ws.Range("H47").Value = 0.4461

$ws.Range("G51").Value = 0.0509
$ws.Range("H51").Value = 0.0509

$ws.Range("G54").Value = 0.098
$ws.Range("H54").Value = 0.098

$ws.Range("G55").Value = 0.0465
$ws.Range("H55").Value = 0.0465

$ws.Range("G63").Value = 1.71
$ws.Range("H63").Value = 1.71

$ws.Range("G64").Value = 0.81236
$ws.Range("H64").Value = 0.81236
